# Update parametric survival model parameter estimates / SEs (and their
# associated covariance-matrix sheets) following a re-run of the individual
# arm model fits and the multivariate NMA data stacking.

$wb = $excel.ActiveWorkbook

# --- Parameter estimate / SE sheets ---------------------------------------

$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.0472551672577
$ws.Range("C2").Value = 0.261938085818681
$ws.Range("B3").Value = 0.155430830791449
$ws.Range("C3").Value = 0.172266877897778

$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.41853795862842
$ws.Range("C2").Value = 0.298670614418002
$ws.Range("B3").Value = -0.997272275177952
$ws.Range("C3").Value = 0.129046366793532

$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.45210266844635
$ws.Range("C2").Value = 0.163509255779092
$ws.Range("B3").Value = 1.76506582427408
$ws.Range("C3").Value = 0.238246846292879

$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -1.9226526959951
$ws.Range("C2").Value = 0.2434248993688
$ws.Range("B3").Value = 0.0148059705014051
$ws.Range("C3").Value = 0.0463504303092982

# "exp" sheet is unchanged.

# --- Covariance matrix sheets ----------------------------------------------

$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0686115608023547
$ws.Range("B2").Value = -0.0339121655313017
$ws.Range("A3").Value = -0.0339121655313017
$ws.Range("B3").Value = 0.0296758772206481

$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.089204135916827
$ws.Range("B2").Value = -0.0321888868121076
$ws.Range("A3").Value = -0.0321888868121076
$ws.Range("B3").Value = 0.0166529647826108

$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0267352767254325
$ws.Range("B2").Value = -0.0101279729218984
$ws.Range("A3").Value = -0.0101279729218984
$ws.Range("B3").Value = 0.0567615597685027

$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0592556816327103
$ws.Range("B2").Value = -0.00869296177512084
$ws.Range("A3").Value = -0.00869296177512084
$ws.Range("B3").Value = 0.00214836238985711

# "exp cov" sheet is unchanged.
